{"js": "// Replace each occurrence of the old problem/date text with the new one.\n// The diff is a set of 1:1 text substitutions (a date label + 25 multiplication\n// problems spread across table cells), each old value is unique in the document,\n// so Body.search + replace is safe and unambiguous.\nconst replacements = [\n  [\"2025-06-19 Thursday\", \"2025-06-20 Friday\"],\n  [\"33\u00d789=\", \"30\u00d758=\"],\n  [\"35\u00d785=\", \"47\u00d735=\"],\n  [\"92\u00d787=\", \"93\u00d786=\"],\n  [\"67\u00d753=\", \"21\u00d769=\"],\n  [\"83\u00d756=\", \"82\u00d776=\"],\n  [\"23\u00d787=\", \"42\u00d729=\"],\n  [\"48\u00d714=\", \"38\u00d783=\"],\n  [\"52\u00d771=\", \"93\u00d769=\"],\n  [\"62\u00d791=\", \"60\u00d781=\"],\n  [\"65\u00d734=\", \"43\u00d738=\"],\n  [\"57\u00d793=\", \"24\u00d777=\"],\n  [\"43\u00d719=\", \"21\u00d748=\"],\n  [\"13\u00d717=\", \"34\u00d778=\"],\n  [\"51\u00d789=\", \"98\u00d732=\"],\n  [\"17\u00d762=\", \"17\u00d788=\"],\n  [\"66\u00d719=\", \"33\u00d756=\"],\n  [\"83\u00d736=\", \"66\u00d714=\"],\n  [\"42\u00d727=\", \"91\u00d740=\"],\n  [\"45\u00d716=\", \"55\u00d738=\"],\n  [\"61\u00d748=\", \"82\u00d722=\"],\n  [\"63\u00d753=\", \"21\u00d714=\"],\n  [\"27\u00d732=\", \"13\u00d755=\"],\n  [\"77\u00d772=\", \"58\u00d719=\"],\n  [\"99\u00d786=\", \"45\u00d743=\"],\n  [\"49\u00d773=\", \"81\u00d740=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each occurrence of the old problem/date text with the new one.\n# The diff is a set of 1:1 text substitutions (a date label + 25 multiplication\n# problems spread across table cells); every old value is unique in the\n# document, so a plain Find/Replace (ReplaceAll) per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = '2025-06-19 Thursday'; New = '2025-06-20 Friday' },\n    @{ Old = '33\u00d789='; New = '30\u00d758=' },\n    @{ Old = '35\u00d785='; New = '47\u00d735=' },\n    @{ Old = '92\u00d787='; New = '93\u00d786=' },\n    @{ Old = '67\u00d753='; New = '21\u00d769=' },\n    @{ Old = '83\u00d756='; New = '82\u00d776=' },\n    @{ Old = '23\u00d787='; New = '42\u00d729=' },\n    @{ Old = '48\u00d714='; New = '38\u00d783=' },\n    @{ Old = '52\u00d771='; New = '93\u00d769=' },\n    @{ Old = '62\u00d791='; New = '60\u00d781=' },\n    @{ Old = '65\u00d734='; New = '43\u00d738=' },\n    @{ Old = '57\u00d793='; New = '24\u00d777=' },\n    @{ Old = '43\u00d719='; New = '21\u00d748=' },\n    @{ Old = '13\u00d717='; New = '34\u00d778=' },\n    @{ Old = '51\u00d789='; New = '98\u00d732=' },\n    @{ Old = '17\u00d762='; New = '17\u00d788=' },\n    @{ Old = '66\u00d719='; New = '33\u00d756=' },\n    @{ Old = '83\u00d736='; New = '66\u00d714=' },\n    @{ Old = '42\u00d727='; New = '91\u00d740=' },\n    @{ Old = '45\u00d716='; New = '55\u00d738=' },\n    @{ Old = '61\u00d748='; New = '82\u00d722=' },\n    @{ Old = '63\u00d753='; New = '21\u00d714=' },\n    @{ Old = '27\u00d732='; New = '13\u00d755=' },\n    @{ Old = '77\u00d772='; New = '58\u00d719=' },\n    @{ Old = '99\u00d786='; New = '45\u00d743=' },\n    @{ Old = '49\u00d773='; New = '81\u00d740=' }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
